$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row text corrections (column labels shifted / typo fixed) ---
$ws.Range("H1").Value = "Indikativ erkennen a)"
$ws.Range("I1").Value = "Indikativ erkennen b)"
$ws.Range("J1").Value = "Indikativ erkennen c)"
$ws.Range("K1").Value = "Indikativ erkennen d)"
$ws.Range("M1").Value = "Indikativ erkennen e)"
$ws.Range("O1").Value = "Konjunktiv II bilden c)"

# --- View state: scrolled right, zoomed in, new active cell ---
$ws.Activate() | Out-Null
$ws.Range("T1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 172
$excel.ActiveWindow.ScrollColumn = 17
$excel.ActiveWindow.ScrollRow = 1
